$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.053.18"
$ws.Range("E2").Value = "  -0.11%  "

# Row 3
$ws.Range("D3").Value = "1.788.53"
$ws.Range("E3").Value = "  -0.04%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "226.77"
$c.ClearFormats()
$ws.Range("E5").Value = "  +1.89%  "

# Row 6
$ws.Range("E6").Value = "  -1.39%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "32.25"
$c.ClearFormats()
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.296"
$c.ClearFormats()
$ws.Range("E9").Value = "  +3.86%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0683"
$c.ClearFormats()

# Row 11
$ws.Range("E11").Value = "  +1.07%  "

# Row 12
$ws.Range("D12").Value = "2.045.09"
$ws.Range("E12").Value = "  -0.07%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "11.35"
$c.ClearFormats()
$ws.Range("E13").Value = "  +3.84%  "

# Row 14
$ws.Range("D14").Value = "1.760.33"
$ws.Range("E14").Value = "  -1.49%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.621"
$c.ClearFormats()
$ws.Range("E15").Value = "  -0.84%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.010.80"
$ws.Range("E16").Value = "  -0.23%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "4.19"
$c.ClearFormats()
$ws.Range("E17").Value = "  +0.37%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "67.83"
$c.ClearFormats()
$ws.Range("E18").Value = "  -0.35%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "242.60"
$c.ClearFormats()
$ws.Range("E19").Value = "  -0.59%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0773"
$ws.Range("E20").Value = "  -1.47%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Range("E21").Value = "  +0.01%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.72"
$c.ClearFormats()
$ws.Range("E22").Value = "  -0.27%  "

# Row 23
$ws.Range("E23").Value = "  -0.11%  "

# Row 24
$ws.Range("E24").Value = "  -2.88%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "162.09"
$c.ClearFormats()
$ws.Range("E25").Value = "  +2.00%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "16.21"
$c.ClearFormats()
$ws.Range("E27").Value = "  -0.93%  "

# Row 28
$ws.Range("E28").Value = "  +0.39%  "

# Row 29
$ws.Range("E29").Value = "  +0.08%  "

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.24"
$c.ClearFormats()
$ws.Range("E30").Value = "  +2.65%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0516"
$c.ClearFormats()
$ws.Range("E31").Value = "  -0.82%  "

# Row 32
$ws.Range("E32").Value = "  -0.78%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.61"
$c.ClearFormats()
$ws.Range("E33").Value = "  +3.14%  "

# Row 34
$ws.Range("E34").Value = "  +1.47%  "

# Row 35
$ws.Range("D35").Value = "1.396.87"

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.653"
$c.ClearFormats()
$ws.Range("E36").Value = "  +0.67%  "

# Row 37
$ws.Range("E37").Value = "  -0.87%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.36"
$c.ClearFormats()
$ws.Range("E38").Value = "  +9.18%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0188"
$c.ClearFormats()
$ws.Range("E39").Value = "  +1.44%  "

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "80.04"
$c.ClearFormats()
$ws.Range("E40").Value = "  +0.49%  "

# Row 41
$ws.Range("E41").Value = "  +0.06%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.920"
$c.ClearFormats()
$ws.Range("E42").Value = "  -0.04%  "

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "13.71"
$c.ClearFormats()
$ws.Range("E43").Value = "  +14.35%  "

# Row 44
$ws.Range("E44").Value = "  -1.10%  "

# Row 45
$ws.Range("E45").Value = "  +8.63%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "6.08"
$c.ClearFormats()
$ws.Range("E46").Value = "  +2.91%  "

# Row 47
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0507"
$c.ClearFormats()
$ws.Range("E47").Value = "  +1.86%  "

# Row 48
$ws.Range("E48").Value = "  +2.65%  "

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "107.63"
$c.ClearFormats()
$ws.Range("E49").Value = "  +0.19%  "

# Row 50
$ws.Range("D50").Value = "1.946.09"
$ws.Range("E50").Value = "  -0.17%  "

# Row 51
$ws.Range("E51").Value = "  -0.11%  "
